# agregamos tamaño a gato
#
# Insert, right after the "Int numeroVidas;" paragraph (and before the
# trailing bookmark paragraph):
#   1. a blank paragraph
#   2. a new "Int tamaño;" paragraph, built with the same
#      proofErr(spellStart)/run/proofErr(spellEnd)/run pattern Word uses
#      for the existing field declarations.

$d = $word.ActiveDocument

# Locate the "Int numeroVidas;" paragraph and the end of its range (i.e.
# right before the paragraph mark that starts the next paragraph). The
# paragraph's Range.Text includes the trailing paragraph-mark (Chr 13),
# so trim it before comparing.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Int numeroVidas;") {
        $target = $p
        break
    }
}

$insertAt = $target.Range.End

# 1) Insert a brand-new, completely empty paragraph right after it.
$r1 = $d.Range($insertAt, $insertAt)
$r1.InsertXML("<blank/>") | Out-Null

# 2) Insert a second empty paragraph right after the first, which will
#    become the new "Int tamaño;" paragraph once populated below.
$insertAt2 = $insertAt + 1
$r2 = $d.Range($insertAt2, $insertAt2)
$r2.InsertXML("<blank/>") | Out-Null

# 3) Fill the second new (still-empty) paragraph with the field
#    declaration runs, matching the existing field-declaration markup.
$fillAt = $insertAt2
$fill = $d.Range($fillAt, $fillAt)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Int</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> tamaño;</w:t></w:r>' +
       '</w:p>'
$fill.InsertXML($xml) | Out-Null
